# Auto-generated Excel COM-interop script
# Applies row-level odds/result corrections for rows 130-145 and 156-163

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ecuador LigaPro Serie A")


# Row 130
$ws.Cells.Item(130, 2).Value = 7483081
$ws.Cells.Item(130, 6).Value = 'Deportivo Cuenca'
$ws.Cells.Item(130, 7).Value = 'El Nacional'
$ws.Cells.Item(130, 8).Value = 1
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 'H'
$ws.Cells.Item(130, 11).Value = 2.75
$ws.Cells.Item(130, 12).Value = 3.25
$ws.Cells.Item(130, 13).Value = 2.55
$ws.Cells.Item(130, 14).Value = 3
$ws.Cells.Item(130, 15).Value = 3.3
$ws.Cells.Item(130, 16).Value = 2.3
$ws.Cells.Item(130, 17).Value = 0.25
$ws.Cells.Item(130, 18).Value = 1.825
$ws.Cells.Item(130, 19).Value = 1.975
$ws.Cells.Item(130, 20).Value = 2.75
$ws.Cells.Item(130, 21).Value = 2
$ws.Cells.Item(130, 22).Value = 1.8
$ws.Cells.Item(130, 23).Value = 2
$ws.Cells.Item(130, 24).Value = -1
$ws.Cells.Item(130, 26).Value = 0.825
$ws.Cells.Item(130, 27).Value = -1
$ws.Cells.Item(130, 28).Value = -1
$ws.Cells.Item(130, 29).Value = 0.8

# Row 131
$ws.Cells.Item(131, 2).Value = 7483281
$ws.Cells.Item(131, 6).Value = 'SD Aucas'
$ws.Cells.Item(131, 7).Value = 'Delfin SC'
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 'D'
$ws.Cells.Item(131, 11).Value = 1.909
$ws.Cells.Item(131, 12).Value = 3.25
$ws.Cells.Item(131, 13).Value = 4.2
$ws.Cells.Item(131, 14).Value = 1.909
$ws.Cells.Item(131, 15).Value = 3.5
$ws.Cells.Item(131, 16).Value = 4
$ws.Cells.Item(131, 17).Value = -0.5
$ws.Cells.Item(131, 18).Value = 1.9
$ws.Cells.Item(131, 19).Value = 1.9
$ws.Cells.Item(131, 21).Value = 1.8
$ws.Cells.Item(131, 22).Value = 2
$ws.Cells.Item(131, 24).Value = 2.5
$ws.Cells.Item(131, 25).Value = -1
$ws.Cells.Item(131, 27).Value = 0.8999999999999999
$ws.Cells.Item(131, 29).Value = 1

# Row 132
$ws.Cells.Item(132, 2).Value = 7483189
$ws.Cells.Item(132, 6).Value = 'Independiente del Valle'
$ws.Cells.Item(132, 7).Value = 'Orense'
$ws.Cells.Item(132, 8).Value = 2
$ws.Cells.Item(132, 9).Value = 2
$ws.Cells.Item(132, 11).Value = 1.4
$ws.Cells.Item(132, 12).Value = 4.75
$ws.Cells.Item(132, 13).Value = 7
$ws.Cells.Item(132, 14).Value = 1.4
$ws.Cells.Item(132, 15).Value = 4.5
$ws.Cells.Item(132, 16).Value = 8
$ws.Cells.Item(132, 17).Value = -1.25
$ws.Cells.Item(132, 18).Value = 1.875
$ws.Cells.Item(132, 19).Value = 1.925
$ws.Cells.Item(132, 21).Value = 1.925
$ws.Cells.Item(132, 22).Value = 1.875
$ws.Cells.Item(132, 24).Value = 3.5
$ws.Cells.Item(132, 27).Value = 0.925
$ws.Cells.Item(132, 28).Value = 0.925
$ws.Cells.Item(132, 29).Value = -1

# Row 133
$ws.Cells.Item(133, 2).Value = 7483247
$ws.Cells.Item(133, 6).Value = 'Mushuc Runa'
$ws.Cells.Item(133, 7).Value = 'Universidad Catolica del Ecuador'
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 2
$ws.Cells.Item(133, 10).Value = 'A'
$ws.Cells.Item(133, 11).Value = 3.25
$ws.Cells.Item(133, 12).Value = 3.2
$ws.Cells.Item(133, 13).Value = 2.25
$ws.Cells.Item(133, 14).Value = 3.5
$ws.Cells.Item(133, 15).Value = 3.25
$ws.Cells.Item(133, 16).Value = 2.1
$ws.Cells.Item(133, 17).Value = 0.5
$ws.Cells.Item(133, 18).Value = 1.775
$ws.Cells.Item(133, 19).Value = 2.025
$ws.Cells.Item(133, 20).Value = 2.5
$ws.Cells.Item(133, 21).Value = 1.9
$ws.Cells.Item(133, 22).Value = 1.9
$ws.Cells.Item(133, 23).Value = -1
$ws.Cells.Item(133, 25).Value = 1.1
$ws.Cells.Item(133, 26).Value = -1
$ws.Cells.Item(133, 27).Value = 1.025
$ws.Cells.Item(133, 29).Value = 0.8999999999999999

# Row 135
$ws.Cells.Item(135, 2).Value = 7482832
$ws.Cells.Item(135, 6).Value = 'Barcelona Guayaquil'
$ws.Cells.Item(135, 7).Value = 'Guayaquil City'
$ws.Cells.Item(135, 8).Value = 2
$ws.Cells.Item(135, 9).Value = 1
$ws.Cells.Item(135, 10).Value = 'H'
$ws.Cells.Item(135, 11).Value = 1.363
$ws.Cells.Item(135, 12).Value = 5
$ws.Cells.Item(135, 13).Value = 7.5
$ws.Cells.Item(135, 14).Value = 1.444
$ws.Cells.Item(135, 15).Value = 4
$ws.Cells.Item(135, 16).Value = 8
$ws.Cells.Item(135, 17).Value = -1.25
$ws.Cells.Item(135, 18).Value = 2.05
$ws.Cells.Item(135, 19).Value = 1.75
$ws.Cells.Item(135, 21).Value = 1.95
$ws.Cells.Item(135, 22).Value = 1.85
$ws.Cells.Item(135, 23).Value = 0.444
$ws.Cells.Item(135, 25).Value = -1
$ws.Cells.Item(135, 26).Value = -0.5
$ws.Cells.Item(135, 27).Value = 0.375
$ws.Cells.Item(135, 28).Value = 0.95
$ws.Cells.Item(135, 29).Value = -1

# Row 136
$ws.Cells.Item(136, 2).Value = 7483188
$ws.Cells.Item(136, 6).Value = 'Gualaceo SC'
$ws.Cells.Item(136, 7).Value = 'Emelec'
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 2
$ws.Cells.Item(136, 10).Value = 'A'
$ws.Cells.Item(136, 11).Value = 3.6
$ws.Cells.Item(136, 12).Value = 3.3
$ws.Cells.Item(136, 13).Value = 2.05
$ws.Cells.Item(136, 14).Value = 2.6
$ws.Cells.Item(136, 15).Value = 3.25
$ws.Cells.Item(136, 16).Value = 2.75
$ws.Cells.Item(136, 17).Value = 0
$ws.Cells.Item(136, 18).Value = 1.8
$ws.Cells.Item(136, 19).Value = 2
$ws.Cells.Item(136, 21).Value = 1.975
$ws.Cells.Item(136, 22).Value = 1.825
$ws.Cells.Item(136, 23).Value = -1
$ws.Cells.Item(136, 25).Value = 1.75
$ws.Cells.Item(136, 26).Value = -1
$ws.Cells.Item(136, 27).Value = 1
$ws.Cells.Item(136, 28).Value = -1
$ws.Cells.Item(136, 29).Value = 0.825

# Row 142
$ws.Cells.Item(142, 2).Value = 7528848
$ws.Cells.Item(142, 6).Value = 'Emelec'
$ws.Cells.Item(142, 7).Value = 'Deportivo Cuenca'
$ws.Cells.Item(142, 8).Value = 2
$ws.Cells.Item(142, 10).Value = 'H'
$ws.Cells.Item(142, 11).Value = 1.75
$ws.Cells.Item(142, 12).Value = 3.5
$ws.Cells.Item(142, 13).Value = 4.2
$ws.Cells.Item(142, 14).Value = 2.4
$ws.Cells.Item(142, 15).Value = 3.1
$ws.Cells.Item(142, 16).Value = 2.75
$ws.Cells.Item(142, 17).Value = -0.25
$ws.Cells.Item(142, 18).Value = 2.05
$ws.Cells.Item(142, 19).Value = 1.75
$ws.Cells.Item(142, 20).Value = 2.25
$ws.Cells.Item(142, 21).Value = 1.8
$ws.Cells.Item(142, 22).Value = 2
$ws.Cells.Item(142, 23).Value = 1.4
$ws.Cells.Item(142, 25).Value = -1
$ws.Cells.Item(142, 26).Value = 1.05
$ws.Cells.Item(142, 27).Value = -1
$ws.Cells.Item(142, 28).Value = 0.8
$ws.Cells.Item(142, 29).Value = -1

# Row 143
$ws.Cells.Item(143, 2).Value = 7528858
$ws.Cells.Item(143, 6).Value = 'Orense'
$ws.Cells.Item(143, 7).Value = 'SD Aucas'
$ws.Cells.Item(143, 8).Value = 1
$ws.Cells.Item(143, 9).Value = 2
$ws.Cells.Item(143, 10).Value = 'A'
$ws.Cells.Item(143, 11).Value = 2.2
$ws.Cells.Item(143, 12).Value = 3.2
$ws.Cells.Item(143, 13).Value = 3.2
$ws.Cells.Item(143, 14).Value = 1.95
$ws.Cells.Item(143, 15).Value = 3.2
$ws.Cells.Item(143, 16).Value = 3.8
$ws.Cells.Item(143, 17).Value = -0.5
$ws.Cells.Item(143, 18).Value = 1.95
$ws.Cells.Item(143, 19).Value = 1.85
$ws.Cells.Item(143, 21).Value = 1.85
$ws.Cells.Item(143, 22).Value = 1.95
$ws.Cells.Item(143, 23).Value = -1
$ws.Cells.Item(143, 25).Value = 2.8
$ws.Cells.Item(143, 26).Value = -1
$ws.Cells.Item(143, 27).Value = 0.8500000000000001
$ws.Cells.Item(143, 28).Value = 0.8500000000000001

# Row 144
$ws.Cells.Item(144, 2).Value = 7528852
$ws.Cells.Item(144, 6).Value = 'Delfin SC'
$ws.Cells.Item(144, 7).Value = 'Tecnico Universitario'
$ws.Cells.Item(144, 8).Value = 2
$ws.Cells.Item(144, 10).Value = 'D'
$ws.Cells.Item(144, 11).Value = 2.1
$ws.Cells.Item(144, 12).Value = 3.4
$ws.Cells.Item(144, 13).Value = 3.1
$ws.Cells.Item(144, 14).Value = 2.1
$ws.Cells.Item(144, 15).Value = 3.4
$ws.Cells.Item(144, 16).Value = 3.1
$ws.Cells.Item(144, 17).Value = -0.25
$ws.Cells.Item(144, 18).Value = 1.8
$ws.Cells.Item(144, 19).Value = 2
$ws.Cells.Item(144, 21).Value = 1.9
$ws.Cells.Item(144, 22).Value = 1.9
$ws.Cells.Item(144, 24).Value = 2.4
$ws.Cells.Item(144, 25).Value = -1
$ws.Cells.Item(144, 26).Value = -0.5
$ws.Cells.Item(144, 27).Value = 0.5
$ws.Cells.Item(144, 28).Value = 0.8999999999999999

# Row 145
$ws.Cells.Item(145, 2).Value = 7528857
$ws.Cells.Item(145, 6).Value = 'Universidad Catolica del Ecuador'
$ws.Cells.Item(145, 7).Value = 'Barcelona Guayaquil'
$ws.Cells.Item(145, 8).Value = 0
$ws.Cells.Item(145, 9).Value = 1
$ws.Cells.Item(145, 10).Value = 'A'
$ws.Cells.Item(145, 11).Value = 1.533
$ws.Cells.Item(145, 12).Value = 4
$ws.Cells.Item(145, 13).Value = 5.5
$ws.Cells.Item(145, 14).Value = 1.5
$ws.Cells.Item(145, 15).Value = 4.333
$ws.Cells.Item(145, 16).Value = 5.25
$ws.Cells.Item(145, 17).Value = -1
$ws.Cells.Item(145, 20).Value = 3
$ws.Cells.Item(145, 21).Value = 1.975
$ws.Cells.Item(145, 22).Value = 1.825
$ws.Cells.Item(145, 24).Value = -1
$ws.Cells.Item(145, 25).Value = 4.25
$ws.Cells.Item(145, 26).Value = -1
$ws.Cells.Item(145, 27).Value = 1
$ws.Cells.Item(145, 28).Value = -1
$ws.Cells.Item(145, 29).Value = 0.825

# Row 156
$ws.Cells.Item(156, 8).Value = 0
$ws.Cells.Item(156, 9).Value = 0
$ws.Cells.Item(156, 10).Value = 'D'
$ws.Cells.Item(156, 14).Value = 1.5
$ws.Cells.Item(156, 15).Value = 4
$ws.Cells.Item(156, 16).Value = 6.5
$ws.Cells.Item(156, 17).Value = -1
$ws.Cells.Item(156, 18).Value = 1.85
$ws.Cells.Item(156, 19).Value = 1.95
$ws.Cells.Item(156, 21).Value = 1.85
$ws.Cells.Item(156, 22).Value = 1.95
$ws.Cells.Item(156, 23).Value = -1
$ws.Cells.Item(156, 24).Value = 3
$ws.Cells.Item(156, 25).Value = -1
$ws.Cells.Item(156, 26).Value = -1
$ws.Cells.Item(156, 27).Value = 0.95
$ws.Cells.Item(156, 28).Value = -1
$ws.Cells.Item(156, 29).Value = 0.95

# Row 157
$ws.Cells.Item(157, 14).Value = 1.444
$ws.Cells.Item(157, 15).Value = 4
$ws.Cells.Item(157, 16).Value = 6
$ws.Cells.Item(157, 17).Value = -1.25
$ws.Cells.Item(157, 18).Value = 2.025
$ws.Cells.Item(157, 19).Value = 1.775
$ws.Cells.Item(157, 21).Value = 1.8
$ws.Cells.Item(157, 22).Value = 2

# Row 158
$ws.Cells.Item(158, 14).Value = 2.5
$ws.Cells.Item(158, 16).Value = 2.625
$ws.Cells.Item(158, 21).Value = 1.85
$ws.Cells.Item(158, 22).Value = 1.95

# Row 159
$ws.Cells.Item(159, 18).Value = 1.775
$ws.Cells.Item(159, 19).Value = 2.025
$ws.Cells.Item(159, 20).Value = 2.25
$ws.Cells.Item(159, 21).Value = 1.8
$ws.Cells.Item(159, 22).Value = 2

# Row 160
$ws.Cells.Item(160, 14).Value = 1.909
$ws.Cells.Item(160, 15).Value = 3.5
$ws.Cells.Item(160, 16).Value = 3.75
$ws.Cells.Item(160, 18).Value = 1.95
$ws.Cells.Item(160, 19).Value = 1.85

# Row 161
$ws.Cells.Item(161, 14).Value = 2
$ws.Cells.Item(161, 15).Value = 3.1
$ws.Cells.Item(161, 16).Value = 3.75
$ws.Cells.Item(161, 18).Value = 2
$ws.Cells.Item(161, 19).Value = 1.8

# Row 162
$ws.Cells.Item(162, 14).Value = 6
$ws.Cells.Item(162, 15).Value = 3.75
$ws.Cells.Item(162, 16).Value = 1.5
$ws.Cells.Item(162, 18).Value = 1.925
$ws.Cells.Item(162, 19).Value = 1.875
$ws.Cells.Item(162, 21).Value = 1.975
$ws.Cells.Item(162, 22).Value = 1.825

# Row 163
$ws.Cells.Item(163, 18).Value = 1.85
$ws.Cells.Item(163, 19).Value = 1.95
